# Weekly refresh of the "Macroferia Regional de Talca - Chirimoya" sheet.
# The underlying records (A,B,C,E,F,G,H,I,J,K,Q,R,T) are unchanged; only the
# per-record Fecha (D), Calidad (L), Volumen (M), Precio minimo/maximo/
# promedio (N/O/P) and Precio $/Kg (S) are refreshed, and one additional
# record (row 31) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2;  D=44455; L="Especial"; M=150; N=30000; O=30000; P=30000; S=3000},
    @{Row=3;  D=44441; L="Primera";  M=150; N=30000; O=30000; P=30000; S=3000},
    @{Row=4;  D=44459; L="Especial"; M=100; N=30000; O=30000; P=30000; S=3000},
    @{Row=5;  D=44462; L="Especial"; M=205; N=30000; O=30000; P=30000; S=3000},
    @{Row=6;  D=44462; L="Primera";  M=180; N=28000; O=28000; P=28000; S=2800},
    @{Row=7;  D=44431; L="Especial"; M=30;  N=30000; O=30000; P=30000; S=3000},
    @{Row=8;  D=44460; L="Especial"; M=80;  N=30000; O=30000; P=30000; S=3000},
    @{Row=9;  D=44446; L="Primera";  M=200; N=30000; O=30000; P=30000; S=3000},
    @{Row=10; D=44463; L="Especial"; M=150; N=30000; O=30000; P=30000; S=3000},
    @{Row=11; D=44463; L="Primera";  M=100; N=26000; O=26000; P=26000; S=2600},
    @{Row=12; D=44447; L="Especial"; M=50;  N=32000; O=32000; P=32000; S=3200},
    @{Row=13; D=44454; L="Especial"; M=320; N=30000; O=30000; P=30000; S=3000},
    @{Row=14; D=44454; L="Primera";  M=300; N=28000; O=28000; P=28000; S=2800},
    @{Row=15; D=44467; L="Especial"; M=100; N=30000; O=30000; P=30000; S=3000},
    @{Row=16; D=44467; L="Primera";  M=100; N=28000; O=28000; P=28000; S=2800},
    @{Row=17; D=44473; L="Primera";  M=200; N=28000; O=28000; P=28000; S=2800},
    @{Row=18; D=44434; L="Especial"; M=60;  N=30000; O=30000; P=30000; S=3000},
    @{Row=19; D=44445; L="Primera";  M=250; N=28000; O=30000; P=29200; S=2920},
    @{Row=20; D=44475; L="Primera";  M=200; N=28000; O=28000; P=28000; S=2800},
    @{Row=21; D=44474; L="Especial"; M=150; N=30000; O=30000; P=30000; S=3000},
    @{Row=22; D=44448; L="Especial"; M=100; N=30000; O=30000; P=30000; S=3000},
    @{Row=23; D=44448; L="Primera";  M=80;  N=28000; O=28000; P=28000; S=2800},
    @{Row=24; D=44453; L="Especial"; M=135; N=30000; O=30000; P=30000; S=3000},
    @{Row=25; D=44435; L="Especial"; M=160; N=30000; O=30000; P=30000; S=3000},
    @{Row=26; D=44468; L="Especial"; M=250; N=30000; O=30000; P=30000; S=3000},
    @{Row=27; D=44466; L="Especial"; M=110; N=30000; O=30000; P=30000; S=3000}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D    # D: Fecha
    $ws.Cells.Item($r, 12).Value = $u.L    # L: Calidad
    $ws.Cells.Item($r, 13).Value = $u.M    # M: Volumen
    $ws.Cells.Item($r, 14).Value = $u.N    # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $u.O    # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $u.P    # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $u.S    # S: Precio $/Kg
}

# New row 31 - a fresh record, copy of the (former) row 2 record's static
# columns with its own date/quality/volume/price data.
$ws.Cells.Item(31, 1).Value  = 5
$ws.Cells.Item(31, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(31, 3).Value  = "Maule"
$ws.Cells.Item(31, 4).Value  = 44432
$ws.Cells.Item(31, 5).Value  = 7
$ws.Cells.Item(31, 6).Value  = "Fruta"
$ws.Cells.Item(31, 7).Value  = 100107
$ws.Cells.Item(31, 8).Value  = "Otros"
$ws.Cells.Item(31, 9).Value  = 100107002
$ws.Cells.Item(31, 10).Value = "Chirimoya"
$ws.Cells.Item(31, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(31, 12).Value = "Especial"
$ws.Cells.Item(31, 13).Value = 70
$ws.Cells.Item(31, 14).Value = 30000
$ws.Cells.Item(31, 15).Value = 30000
$ws.Cells.Item(31, 16).Value = 30000
$ws.Cells.Item(31, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(31, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 19).Value = 3000
$ws.Cells.Item(31, 20).Value = 10

# Match the date formatting/style used by the other Fecha cells (D2:D30).
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(30, 4).NumberFormat
